$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 49.88947433333334
$ws.Range("H2").Value = 149.668423
$ws.Range("I2").Value = 0.2324880572195875
$ws.Range("J2").Value = 0.2324880572195874
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5034623333333333
$ws.Range("N2").Value = 1.510387
$ws.Range("O2").Value = 0.3640009504979511
$ws.Range("P2").Value = 0.364000950497951
$ws.Range("Q2").Value = 25.11747115663345
$ws.Range("R2").Value = 226.0572404097011
$ws.Range("S2").Value = 0.08462587380735187
$ws.Range("T2").Value = 0.08462587380735184

# Row 3
$ws.Range("G3").Value = 49.88947433333334
$ws.Range("H3").Value = 149.668423
$ws.Range("I3").Value = 0.2324880572195875
$ws.Range("J3").Value = 0.2324880572195874
$ws.Range("O3").Value = 0.04322018294675573
$ws.Range("P3").Value = 0.04322018294675572
$ws.Range("Q3").Value = 2.982359515997111
$ws.Range("R3").Value = 26.841235643974
$ws.Range("S3").Value = 0.01004817636596638
$ws.Range("T3").Value = 0.01004817636596638

# Row 4
$ws.Range("G4").Value = 49.88947433333334
$ws.Range("H4").Value = 149.668423
$ws.Range("I4").Value = 0.2324880572195875
$ws.Range("J4").Value = 0.2324880572195874
$ws.Range("M4").Value = 0.819893
$ws.Range("N4").Value = 2.459679
$ws.Range("O4").Value = 0.5927788665552932
$ws.Range("P4").Value = 0.5927788665552931
$ws.Range("Q4").Value = 40.90403077957967
$ws.Range("R4").Value = 368.136277016217
$ws.Range("S4").Value = 0.1378140070462692
$ws.Range("T4").Value = 0.1378140070462692

# Row 5
$ws.Range("I5").Value = 0.295249080025651
$ws.Range("J5").Value = 0.295249080025651
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5034623333333333
$ws.Range("N5").Value = 1.510387
$ws.Range("O5").Value = 0.3640009504979511
$ws.Range("P5").Value = 0.364000950497951
$ws.Range("Q5").Value = 31.89802667825834
$ws.Range("R5").Value = 287.082240104325
$ws.Range("S5").Value = 0.1074709457629826
$ws.Range("T5").Value = 0.1074709457629826

# Row 6
$ws.Range("I6").Value = 0.295249080025651
$ws.Range("J6").Value = 0.295249080025651
$ws.Range("O6").Value = 0.04322018294675573
$ws.Range("P6").Value = 0.04322018294675572
$ws.Range("S6").Value = 0.01276071925356996
$ws.Range("T6").Value = 0.01276071925356996

# Row 7
$ws.Range("I7").Value = 0.295249080025651
$ws.Range("J7").Value = 0.295249080025651
$ws.Range("M7").Value = 0.819893
$ws.Range("N7").Value = 2.459679
$ws.Range("O7").Value = 0.5927788665552932
$ws.Range("P7").Value = 0.5927788665552931
$ws.Range("Q7").Value = 51.946227266225
$ws.Range("R7").Value = 467.516045396025
$ws.Range("S7").Value = 0.1750174150090985
$ws.Range("T7").Value = 0.1750174150090985

# Row 8
$ws.Range("G8").Value = 52.37451933333333
$ws.Range("H8").Value = 157.123558
$ws.Range("I8").Value = 0.2440685216737345
$ws.Range("J8").Value = 0.2440685216737345
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5034623333333333
$ws.Range("N8").Value = 1.510387
$ws.Range("O8").Value = 0.3640009504979511
$ws.Range("P8").Value = 0.364000950497951
$ws.Range("Q8").Value = 26.36859771077178
$ws.Range("R8").Value = 237.317379396946
$ws.Range("S8").Value = 0.08884117387586912
$ws.Range("T8").Value = 0.08884117387586911

# Row 9
$ws.Range("G9").Value = 52.37451933333333
$ws.Range("H9").Value = 157.123558
$ws.Range("I9").Value = 0.2440685216737345
$ws.Range("J9").Value = 0.2440685216737345
$ws.Range("O9").Value = 0.04322018294675573
$ws.Range("P9").Value = 0.04322018294675572
$ws.Range("Q9").Value = 3.130913849400444
$ws.Range("R9").Value = 28.178224644604
$ws.Range("S9").Value = 0.01054868615828302
$ws.Range("T9").Value = 0.01054868615828302

# Row 10
$ws.Range("G10").Value = 52.37451933333333
$ws.Range("H10").Value = 157.123558
$ws.Range("I10").Value = 0.2440685216737345
$ws.Range("J10").Value = 0.2440685216737345
$ws.Range("M10").Value = 0.819893
$ws.Range("N10").Value = 2.459679
$ws.Range("O10").Value = 0.5927788665552932
$ws.Range("P10").Value = 0.5927788665552931
$ws.Range("Q10").Value = 42.94150177976466
$ws.Range("R10").Value = 386.473516017882
$ws.Range("S10").Value = 0.1446786616395823
$ws.Range("T10").Value = 0.1446786616395823

# Row 11
$ws.Range("G11").Value = 48.96808833333333
$ws.Range("H11").Value = 146.904265
$ws.Range("I11").Value = 0.2281943410810271
$ws.Range("J11").Value = 0.228194341081027
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.5034623333333333
$ws.Range("N11").Value = 1.510387
$ws.Range("O11").Value = 0.3640009504979511
$ws.Range("P11").Value = 0.364000950497951
$ws.Range("Q11").Value = 24.65358801117278
$ws.Range("R11").Value = 221.882292100555
$ws.Range("S11").Value = 0.08306295705174749
$ws.Range("T11").Value = 0.08306295705174747

# Row 12
$ws.Range("G12").Value = 48.96808833333333
$ws.Range("H12").Value = 146.904265
$ws.Range("I12").Value = 0.2281943410810271
$ws.Range("J12").Value = 0.228194341081027
$ws.Range("O12").Value = 0.04322018294675573
$ws.Range("P12").Value = 0.04322018294675572
$ws.Range("Q12").Value = 2.927279675174444
$ws.Range("R12").Value = 26.34551707657
$ws.Range("S12").Value = 0.009862601168936365
$ws.Range("T12").Value = 0.009862601168936364

# Row 13
$ws.Range("G13").Value = 48.96808833333333
$ws.Range("H13").Value = 146.904265
$ws.Range("I13").Value = 0.2281943410810271
$ws.Range("J13").Value = 0.228194341081027
$ws.Range("M13").Value = 0.819893
$ws.Range("N13").Value = 2.459679
$ws.Range("O13").Value = 0.5927788665552932
$ws.Range("P13").Value = 0.5927788665552931
$ws.Range("Q13").Value = 40.14859284788167
$ws.Range("R13").Value = 361.337335630935
$ws.Range("S13").Value = 0.1352687828603432
$ws.Range("T13").Value = 0.1352687828603432
